# Update gh-pages output data (想去人数 / F column) for 南宁-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 3-9 in column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 307
$wsExpo.Range("F4").Value = 212
$wsExpo.Range("F5").Value = 2621
$wsExpo.Range("F6").Value = 1873
$wsExpo.Range("F7").Value = 359
$wsExpo.Range("F8").Value = 112
$wsExpo.Range("F9").Value = 926

# Sheet "全部类型" (sheet4): rows 3-7, 9-10 in column F (row 8 unrelated, unchanged)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 307
$wsAll.Range("F4").Value = 212
$wsAll.Range("F5").Value = 2621
$wsAll.Range("F6").Value = 1873
$wsAll.Range("F7").Value = 359
$wsAll.Range("F9").Value = 112
$wsAll.Range("F10").Value = 926
